# C5-PowerPoint.pptx edit
#
# Target diff (summarised):
#   1. ppt/slides/slide6.xml - the "SOURCES OF FINANCE" table's
#      <a:tableStyleId> is switched from the custom Google-Slides style
#      {2597ECAA-AA4C-4DD6-8202-859E3E2AE77C} to the built-in PowerPoint
#      table style {680DE56F-00D5-4ADE-BC47-AEE41A4FEB7B}.
#   2. ppt/presentation.xml loses the embedTrueTypeFonts="1" attribute and
#      the <p:embeddedFontLst> block (the embedded "Limelight" font stops
#      being flagged for embedding).
#   3. ppt/theme/theme1.xml and ppt/theme/theme2.xml swap their contents
#      (which theme part number holds "Office Theme" vs "Integral").
#
# Items 2 and 3 are low-level OOXML-packaging side effects (PowerPoint's
# own save "churn"/font-embedding toggle) that are not reachable through
# the Presentation/Application COM object model - PowerPoint does not
# expose Presentation.EmbedTrueTypeFonts, a way to drop
# <p:embeddedFontLst>, or any API to rewrite/renumber raw theme parts.
# (Font.Embeddable/Font.Embedded on a TextRange's Font are themselves
# read-only status flags, not switches.) Only the table style change
# below is something a user/automation can actually drive via the object
# model, so that is what this script performs.

$p = $ppt.ActivePresentation

# The table lives on slide 6 ("SOURCES OF FINANCE"), as the 2nd shape
# (1: title placeholder, 2: graphicFrame holding the table).
$targetStyleId = "{680DE56F-00D5-4ADE-BC47-AEE41A4FEB7B}"
$oldStyleId = "{2597ECAA-AA4C-4DD6-8202-859E3E2AE77C}"

$applied = $false
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style -eq $oldStyleId) {
                $table.ApplyStyle($targetStyleId)
                Write-Host "Slide $i shape $j ($($shape.Name)): table style ->" $table.Style
                $applied = $true
            }
        }
    }
}

if (-not $applied) {
    Write-Host "WARNING: target table (style $oldStyleId) was not found."
}
